$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(5)
$tr = $shp.TextFrame.TextRange

# --- Paragraph 2: "Not published online, only on localhost" ---
# Add trailing space to the existing run, then append a differently
# formatted run with the new sentence.
$para2 = $tr.Paragraphs(2)
$r1 = $para2.Runs(1)
$oldR1Text = "Not published online, only on localhost "
$r1.Text = $oldR1Text

$startOfPara2 = $para2.Start
$newText2 = "and an error while displaying forecasted output."
$para2.InsertAfter($newText2) | Out-Null

$newRun2 = $tr.Characters($startOfPara2 + $oldR1Text.Length, $newText2.Length)
$newRun2.Font.Size = 18
$newRun2.Font.Bold = $false
$newRun2.Font.Italic = $false
$newRun2.Font.Underline = $false
$newRun2.Font.Strike = 0
$newRun2.Font.Shadow = $false
$newRun2.Font.Color.RGB = 0
$newRun2.Font.Name = "Arial"

# --- Paragraph 4: "Change enter latitude, longitude, ..." ---
# Collapse the three runs into a single run with uniform text/formatting.
$para4 = $tr.Paragraphs(4)
$startOfPara4 = $para4.Start
$run1Text = "Change enter latitude, longitude, to select on map "
$run2Text = "which picks latitude and longitude as "
$run3Text = "input"

$tailRange = $tr.Characters($startOfPara4 + $run1Text.Length, $run2Text.Length + $run3Text.Length)
$tailRange.Delete()

$para4Run1 = $para4.Runs(1)
$para4Run1.Text = $run1Text + $run2Text + $run3Text

# --- Resize the textbox to fit the new content ---
$shp.Height = 1107996 / 12700
